# Quiz workbook update:
#   - turn the Question/Answer range into a proper Excel table (Table1)
#   - re-sort the quiz alphabetically by Question (this is what brings the
#     "svensk farmor" / potato-hat riddle, already present in the sheet,
#     into its alphabetically-sorted spot) so every answer can be browsed
#     via the table's filter/sort UI
#   - nudge column B's width and leave the selection where the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn A1:B20 (Question/Answer, with header row) into a real table.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:B20"), $null, 1)
$tbl.Name = "Table1"

# Sort the table alphabetically by the Question column.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("A2:A20"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Widen column B a bit, matching the author's manual resize.
$ws.Columns(2).ColumnWidth = 8.25

# Restore the last-used selection.
[void]$ws.Range("F9").Select()
